$d = $word.ActiveDocument

# 1. Title heading + bold byline near end (Replace:=wdReplaceAll replaces every occurrence in one call)
$range = $d.Content
$range.Find.Execute("Play La Mafia Heist Free: Review & Gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Play La Mafia Heist and Claim Exciting Free Spins - Review", 2)

# 2. "What we like" bullets
$range = $d.Content
$range.Find.Execute("Impressive graphics and well-detailed characters", $true, $false, $false, $false, $false, $true, 1, $false, "Impressive comic book-style graphics", 2)

$range = $d.Content
$range.Find.Execute("Two exciting bonus features", $true, $false, $false, $false, $false, $true, 1, $false, "Two exciting bonus features - free spins and Gold Reserve", 2)

$range = $d.Content
$range.Find.Execute("High volatility perfect for high-risk players", $true, $false, $false, $false, $false, $true, 1, $false, "High volatility for high-risk players", 2)

$range = $d.Content
$range.Find.Execute("Autoplay with up to 500 continuous spins available", $true, $false, $false, $false, $false, $true, 1, $false, "Autoplay function for convenience", 2)

# 3. "What we don't like" bullets
$range = $d.Content
$range.Find.Execute("Limited bet options for low rollers", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting options with a minimum bet of €0.20", 2)

$range = $d.Content
$range.Find.Execute("No progressive jackpot", $true, $false, $false, $false, $false, $true, 1, $false, "Only 40 fixed paylines", 2)

# 4. Meta description sentence
$range = $d.Content
$range.Find.Execute("Read our review of La Mafia Heist, a high-volatility 5x3 slot game with 40 fixed paylines. Play for free and enjoy two exciting bonus features.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of La Mafia Heist, a slot game with impressive graphics and exciting free spins. Play now for free!", 2)
